$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily "Diaria" data rows appended after the existing last row (208),
# one row per new business day reported by MV (see commit message).
# Columns: A = Serie (date label, stored as text), C = 2 años, D = 5 años, E = 10 años.
$newRows = @(
    @{ Row = 209; Date = "02-11-2021"; C = $null; D = 3.58; E = 3.43 },
    @{ Row = 210; Date = "03-11-2021"; C = 5.26;  D = 3.59; E = 3.43 },
    @{ Row = 211; Date = "04-11-2021"; C = $null; D = 3.6;  E = 3.5  },
    @{ Row = 212; Date = "05-11-2021"; C = 5.56;  D = 3.67; E = 3.5  }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Write the date label as a literal text value (matching the existing
    # "dd-mm-yyyy" text entries in column A) without letting Excel's
    # autodetection turn it into a date serial number / date-formatted cell:
    # compute it via a formula that yields a text result, then freeze that
    # computed text in place as a plain value (Copy + PasteSpecial values),
    # which preserves the default (unstyled) cell format.
    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.Formula = '="' + $r.Date + '"'
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)  # xlPasteValues

    if ($null -ne $r.C) {
        $ws.Cells.Item($rowNum, 3).Value = $r.C
    }
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
}

$excel.CutCopyMode = 0
